$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update to candidate id=2's scoring data ---
$ws.Range("B2").Value2 = 2
$ws.Range("C2").Value2 = 56.54
$ws.Range("D2").Value2 = 0.85
$e2Text = @'


Reasoning: The candidate has the required skills for the job, but the projects they have worked on do not demonstrate any experience in NLP, Pytorch, Computer Vision, or Python. The candidate does have experience with ReactJS, HTML, CSS, Flutter, Dart, and Firebase, however, these are not the skills required for the job. Therefore, the candidate is only partially suited for the role and has been given a score of 56.54.
'@
$ws.Range("E2").Value = $e2Text

# --- Row 3: update to candidate id=3's scoring data ---
$ws.Range("B3").Value2 = 3
$ws.Range("C3").Value2 = 81.56
$ws.Range("D3").Value2 = 1
$e3Text = @'


The candidate has been given a score of 81.56 because they have the technical skills and experience in the required fields necessary for the NLP engineer position, such as NLP, Pytorch, Computer Vision, and Python. Their projects demonstrate their knowledge in data analysis, multi-model data analysis, object detection, text recognition, and web development. This makes them an excellent fit for the job and they have been given a high score accordingly.
'@
$ws.Range("E3").Value = $e3Text

# --- Remove old row 4 (its candidate is now folded into row 3 above) ---
$ws.Rows(4).Delete()
